$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1 (Post-conditions bullet):
#   "Incident state changed to Closed"
#     -> "Incident state changed to " (run 1) + "Completed" (run 2)
# ---------------------------------------------------------------------------
$rng1 = $d.Content.Duplicate
$found1 = $rng1.Find.Execute("Incident state changed to Closed", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Incident state changed to Closed'"
}
# Collapse onto just the trailing "Closed" word that needs to become "Completed".
$closedRng = $rng1.Duplicate
$closedRng.MoveStart(1, $rng1.End - $rng1.Start - 6) | Out-Null   # wdCharacter = 1; "Closed" is 6 chars

# Toggling Bold on/off around the text assignment forces the edited text to
# stay in its own run (instead of being re-merged with its neighbor) while
# dropping the rsid attributes from this freshly-edited run -- exactly how
# Word marks newly authored content that was typed in this editing session.
$closedRng.Font.Bold = 1
$closedRng.Text = "Completed"
$closedRng.Font.Bold = 0

# ---------------------------------------------------------------------------
# Edit 2 (Main flow step, "System changes the status..."):
#   "System changes the status of incident to 'closed'."
#     -> "System changes the status of incident to '" (run 1)
#        + "Completed" (run 2) + "'" (run 3) + "." (run 4)
# ---------------------------------------------------------------------------
$prefix = $d.Content.Duplicate
$found2 = $prefix.Find.Execute("System changes the status of incident to " + [char]0x2018, $true, $false, `
    $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find 'System changes the status of incident to lsquo'"
}

# "closed" immediately follows the opening curly quote we just matched.
$rClosed = $d.Range($prefix.End, $prefix.End + 6)
$rClosed.Font.Bold = 1
$rClosed.Text = "Completed"
$rClosed.Font.Bold = 0

# Closing curly quote right after "Completed" -- split into its own run.
$quoteStart = $rClosed.End
$rQuote = $d.Range($quoteStart, $quoteStart + 1)
$rQuote.Font.Bold = 1
$rQuote.Text = "X"
$rQuote.Text = [char]0x2019
$rQuote.Font.Bold = 0

# Trailing period -- split into its own run too.
$periodStart = $rQuote.End
$rPeriod = $d.Range($periodStart, $periodStart + 1)
$rPeriod.Font.Bold = 1
$rPeriod.Text = "X"
$rPeriod.Text = "."
$rPeriod.Font.Bold = 0
